$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM row "R7 R8 R9 R10 " (qty 4, RES 100K OHM 1/8W 3.6MM) is being split
# into two rows: "R7 R8 " keeps the 100K value, "R9 R10 " becomes 6.8K OHM.

# 1) Insert a new row below the existing R7/R8/R9/R10 row (row 22) so the rest
#    of the table (RV1, RV2 RV3, SW1, U1, U2 U3, TUBE, HEATSINK, KNOB) shifts
#    down by one row, and the new row inherits row 21's formatting.
$ws.Rows.Item(22).Insert()

# 2) Row 21 becomes the "R9 R10" line with the new 6.8K value (link unchanged).
$ws.Range("A21").Value = "R9 R10 "
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "RES 6.8K OHM 1/8W 3.6MM"

# 3) Row 22 (the newly inserted row) becomes the "R7 R8" line, keeping the
#    original 100K value and hyperlink target.
$ws.Range("A22").Value = "R7 R8 "
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = "RES 100K OHM 1/8W 3.6MM"
$ws.Range("D22").Value = "https://www.aliexpress.com/item/33007959640.html"
$ws.Hyperlinks.Add($ws.Range("D22"), "https://www.aliexpress.com/item/33007959640.html")

# 4) Match the recorded selection state left behind in the saved file.
$ws.Range("E:I").Select()
